# Rename the worksheet (was "RGossF-HW10.xpc")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "RGossF"

# Tiny floating point corrections to existing cells
$ws.Range("I13").Value = 0.9933815826730794
$ws.Range("H15").Value = 0.9955707065057025

# Append a new data row (row 16) for the Gaussian Quadrature scheme result
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.636273015090748
$ws.Range("D16").Value = 2.044356414943652
$ws.Range("E16").Value = 1.021648733139336
$ws.Range("F16").Value = 1.636273015090748
$ws.Range("G16").Value = 0.7430223234523844
$ws.Range("H16").Value = 2.018785056086185
$ws.Range("I16").Value = 0.7717809265185523
$ws.Range("J16").Value = 2.044356414943652
$ws.Range("K16").Value = 1.533002574041494
$ws.Range("L16").Value = 1.584637794566121
$ws.Range("M16").Value = 1.372644411538476

# Match formatting of the row above (border/bold/center style) for the new index cell
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
